$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 98654
$ws.Range("E2").Value = 10441
$ws.Range("F2").Value = 10441
$ws.Range("G2").Value = 9557
$ws.Range("H2").Value = 8343
$ws.Range("I2").Value = 8445
$ws.Range("J2").Value = -103
$ws.Range("K2").Value = 173327
$ws.Range("L2").Value = 47469
$ws.Range("M2").Value = 125858
$ws.Range("N2").Value = 122519
$ws.Range("O2").Value = 3339
$ws.Range("P2").Value = 8794
$ws.Range("Q2").Value = 6003
$ws.Range("R2").Value = -7209
$ws.Range("S2").Value = -599
$ws.Range("T2").Value = 3608
$ws.Range("U2").Value = 2395
$ws.Range("V2").Value = 18774
$ws.Range("W2").Value = 10.58
$ws.Range("X2").Value = 8.46
$ws.Range("Y2").Value = 7.03
$ws.Range("Z2").Value = 4.92
$ws.Range("AA2").Value = 37.72
$ws.Range("AB2").Value = 1309.58
$ws.Range("AC2").Value = 4802
$ws.Range("AD2").Value = 12.74
$ws.Range("AE2").Value = 69704
$ws.Range("AF2").Value = 0.88
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1.63
$ws.Range("AI2").Value = 20.83
$ws.Range("AJ2").Value = 172557131

$ws.Range("D3").Value = 99654
$ws.Range("E3").Value = 11386
$ws.Range("F3").Value = 11380
$ws.Range("G3").Value = 10761
$ws.Range("H3").Value = 9438
$ws.Range("I3").Value = 9442
$ws.Range("J3").Value = -4
$ws.Range("K3").Value = 181234
$ws.Range("L3").Value = 48152
$ws.Range("M3").Value = 133082
$ws.Range("N3").Value = 129755
$ws.Range("O3").Value = 3327
$ws.Range("P3").Value = 8794
$ws.Range("Q3").Value = 8863
$ws.Range("R3").Value = -3683
$ws.Range("S3").Value = -1452
$ws.Range("T3").Value = 2630
$ws.Range("U3").Value = 6233
$ws.Range("V3").Value = 19165
$ws.Range("W3").Value = 11.43
$ws.Range("X3").Value = 9.47
$ws.Range("Y3").Value = 7.49
$ws.Range("Z3").Value = 5.32
$ws.Range("AA3").Value = 36.18
$ws.Range("AB3").Value = 1391.17
$ws.Range("AC3").Value = 5369
$ws.Range("AD3").Value = 13.19
$ws.Range("AE3").Value = 73820
$ws.Range("AF3").Value = 0.96
$ws.Range("AG3").Value = 1300
$ws.Range("AH3").Value = 1.84
$ws.Range("AI3").Value = 24.22
$ws.Range("AJ3").Value = 172557131

$ws.Range("D4").Value = 96739
$ws.Range("E4").Value = 13227
$ws.Range("F4").Value = 13447
$ws.Range("G4").Value = 12905
$ws.Range("H4").Value = 10913
$ws.Range("I4").Value = 10748
$ws.Range("J4").Value = 165
$ws.Range("K4").Value = 194894
$ws.Range("L4").Value = 52694
$ws.Range("M4").Value = 142201
$ws.Range("N4").Value = 138744
$ws.Range("O4").Value = 3457
$ws.Range("P4").Value = 8794
$ws.Range("Q4").Value = 8782
$ws.Range("R4").Value = -3438
$ws.Range("S4").Value = -2754
$ws.Range("T4").Value = 2745
$ws.Range("U4").Value = 6037
$ws.Range("V4").Value = 18776
$ws.Range("W4").Value = 13.67
$ws.Range("X4").Value = 11.28
$ws.Range("Y4").Value = 8.01
$ws.Range("Z4").Value = 5.8
$ws.Range("AA4").Value = 37.06
$ws.Range("AB4").Value = 1496.2
$ws.Range("AC4").Value = 6111
$ws.Range("AD4").Value = 9.82
$ws.Range("AE4").Value = 78934
$ws.Range("AF4").Value = 0.76
$ws.Range("AG4").Value = 1300
$ws.Range("AH4").Value = 2.17
$ws.Range("AI4").Value = 21.28
$ws.Range("AJ4").Value = 172557131

$ws.Range("D5").Value = 118411
$ws.Range("E5").Value = 21858
$ws.Range("F5").Value = 21858
$ws.Range("G5").Value = 27413
$ws.Range("H5").Value = 24356
$ws.Range("I5").Value = 23959
$ws.Range("J5").Value = 397
$ws.Range("K5").Value = 216477
$ws.Range("L5").Value = 54682
$ws.Range("M5").Value = 161795
$ws.Range("N5").Value = 160026
$ws.Range("O5").Value = 1770
$ws.Range("P5").Value = 8794
$ws.Range("Q5").Value = 10584
$ws.Range("R5").Value = -6893
$ws.Range("S5").Value = -1398
$ws.Range("T5").Value = 2921
$ws.Range("U5").Value = 7663
$ws.Range("V5").Value = 13927
$ws.Range("W5").Value = 18.46
$ws.Range("X5").Value = 20.57
$ws.Range("Y5").Value = 16.04
$ws.Range("Z5").Value = 11.84
$ws.Range("AA5").Value = 33.8
$ws.Range("AB5").Value = 1748.54
$ws.Range("AC5").Value = 13623
$ws.Range("AD5").Value = 6.68
$ws.Range("AE5").Value = 91042
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 1300
$ws.Range("AH5").Value = 1.43
$ws.Range("AI5").Value = 9.54
$ws.Range("AJ5").Value = 172557131

$ws.Range("D6").Value = 119448
$ws.Range("E6").Value = 19638
$ws.Range("F6").Value = 19638
$ws.Range("G6").Value = 19301
$ws.Range("H6").Value = 18829
$ws.Range("I6").Value = 18639
$ws.Range("K6").Value = 228359
$ws.Range("L6").Value = 46714
$ws.Range("M6").Value = 181645
$ws.Range("N6").Value = 179987
$ws.Range("P6").Value = 8794
$ws.Range("Q6").Value = 4949
$ws.Range("R6").Value = -6586
$ws.Range("S6").Value = -1264
$ws.Range("T6").Value = 2635
$ws.Range("U6").Value = 2314
$ws.Range("V6").Value = 15195
$ws.Range("W6").Value = 16.44
$ws.Range("X6").Value = 15.76
$ws.Range("Y6").Value = 10.96
$ws.Range("Z6").Value = 8.46
$ws.Range("AA6").Value = 25.72
$ws.Range("AB6").Value = 1980.32
$ws.Range("AC6").Value = 10598
$ws.Range("AD6").Value = 6.6
$ws.Range("AE6").Value = 102399
$ws.Range("AF6").Value = 0.68
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 2.86
$ws.Range("AI6").Value = 18.87
$ws.Range("AJ6").Value = 172557131

$ws.Range("D7").Value = 69770
$ws.Range("E7").Value = 14702
$ws.Range("G7").Value = 14449
$ws.Range("H7").Value = 16512
$ws.Range("I7").Value = 16112
$ws.Range("K7").Value = 229619
$ws.Range("L7").Value = 33975
$ws.Range("M7").Value = 195644
$ws.Range("N7").Value = 193634
$ws.Range("P7").Value = 8791
$ws.Range("Q7").Value = 12110
$ws.Range("R7").Value = 4978
$ws.Range("S7").Value = -2865
$ws.Range("T7").Value = 444
$ws.Range("U7").Value = 14608
$ws.Range("W7").Value = 21.07
$ws.Range("X7").Value = 23.67
$ws.Range("Y7").Value = 8.63
$ws.Range("Z7").Value = 7.21
$ws.Range("AA7").Value = 17.37
$ws.Range("AC7").Value = 9161
$ws.Range("AD7").Value = 7.67
$ws.Range("AE7").Value = 110162
$ws.Range("AF7").Value = 0.64
$ws.Range("AG7").Value = 2188
$ws.Range("AH7").Value = 3.11
$ws.Range("AI7").Value = 23.43

$ws.Range("D8").Value = 75841
$ws.Range("E8").Value = 18518
$ws.Range("G8").Value = 18403
$ws.Range("H8").Value = 16991
$ws.Range("I8").Value = 16771
$ws.Range("K8").Value = 242979
$ws.Range("L8").Value = 33744
$ws.Range("M8").Value = 209235
$ws.Range("N8").Value = 207137
$ws.Range("P8").Value = 8791
$ws.Range("Q8").Value = 15267
$ws.Range("R8").Value = -3416
$ws.Range("S8").Value = -4069
$ws.Range("T8").Value = 875
$ws.Range("U8").Value = 16700
$ws.Range("W8").Value = 24.42
$ws.Range("X8").Value = 22.4
$ws.Range("Y8").Value = 8.37
$ws.Range("Z8").Value = 7.19
$ws.Range("AA8").Value = 16.13
$ws.Range("AC8").Value = 9536
$ws.Range("AD8").Value = 7.37
$ws.Range("AE8").Value = 117845
$ws.Range("AF8").Value = 0.6
$ws.Range("AG8").Value = 2400
$ws.Range("AH8").Value = 3.41
$ws.Range("AI8").Value = 24.69

$ws.Range("D9").Value = 81398
$ws.Range("E9").Value = 20913
$ws.Range("G9").Value = 20708
$ws.Range("H9").Value = 18763
$ws.Range("I9").Value = 18511
$ws.Range("K9").Value = 256369
$ws.Range("L9").Value = 32022
$ws.Range("M9").Value = 224348
$ws.Range("N9").Value = 222218
$ws.Range("P9").Value = 8791
$ws.Range("Q9").Value = 16620
$ws.Range("R9").Value = -3810
$ws.Range("S9").Value = -4246
$ws.Range("T9").Value = 875
$ws.Range("U9").Value = 18895
$ws.Range("W9").Value = 25.69
$ws.Range("X9").Value = 23.05
$ws.Range("Y9").Value = 8.62
$ws.Range("Z9").Value = 7.51
$ws.Range("AA9").Value = 14.27
$ws.Range("AC9").Value = 10525
$ws.Range("AD9").Value = 6.68
$ws.Range("AE9").Value = 126425
$ws.Range("AF9").Value = 0.56
$ws.Range("AG9").Value = 2500
$ws.Range("AH9").Value = 3.56
$ws.Range("AI9").Value = 23.3
